$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.040764356737687
$ws.Cells.Item(2, 4).Value = 1.04851290669562
$ws.Cells.Item(2, 5).Value = 1.049271138652218
$ws.Cells.Item(2, 6).Value = 1.060499108654025
$ws.Cells.Item(2, 9).Value = 1.038760242541983
$ws.Cells.Item(2, 10).Value = 1.045849200829952
$ws.Cells.Item(2, 11).Value = 1.051272543180697
$ws.Cells.Item(2, 12).Value = 1.052028661921875
$ws.Cells.Item(2, 13).Value = 1.063225722269284
$ws.Cells.Item(2, 14).Value = 1.04733442665571
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.041734690530341
$ws.Cells.Item(3, 4).Value = 1.049260494261625
$ws.Cells.Item(3, 5).Value = 1.050112676945145
$ws.Cells.Item(3, 6).Value = 1.061374458915873
$ws.Cells.Item(3, 9).Value = 1.038927629996514
$ws.Cells.Item(3, 10).Value = 1.046465182399192
$ws.Cells.Item(3, 11).Value = 1.05183235797746
$ws.Cells.Item(3, 12).Value = 1.052682337049213
$ws.Cells.Item(3, 13).Value = 1.063915356078076
$ws.Cells.Item(3, 14).Value = 1.047951282989433
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.042363190272252
$ws.Cells.Item(4, 4).Value = 1.049744596310185
$ws.Cells.Item(4, 5).Value = 1.050658074965333
$ws.Cells.Item(4, 6).Value = 1.061941640165175
$ws.Cells.Item(4, 9).Value = 1.039034774585782
$ws.Cells.Item(4, 10).Value = 1.046863775842855
$ws.Cells.Item(4, 11).Value = 1.052194292140687
$ws.Cells.Item(4, 12).Value = 1.053105527441863
$ws.Cells.Item(4, 13).Value = 1.064361724062979
$ws.Cells.Item(4, 14).Value = 1.048350442481505
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.042627561086979
$ws.Cells.Item(5, 4).Value = 1.049948198209434
$ws.Cells.Item(5, 5).Value = 1.050887566275619
$ws.Cells.Item(5, 6).Value = 1.06218026621159
$ws.Cells.Item(5, 9).Value = 1.039079538388775
$ws.Cells.Item(5, 10).Value = 1.047031346395825
$ws.Cells.Item(5, 11).Value = 1.052346375486543
$ws.Cells.Item(5, 12).Value = 1.0532834879219
$ws.Cells.Item(5, 13).Value = 1.064549406700378
$ws.Cells.Item(5, 14).Value = 1.04851825100388
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.042671958848197
$ws.Cells.Item(6, 4).Value = 1.049982388828908
$ws.Cells.Item(6, 5).Value = 1.050926110913126
$ws.Cells.Item(6, 6).Value = 1.062220343263752
$ws.Cells.Item(6, 9).Value = 1.039087038005129
$ws.Cells.Item(6, 10).Value = 1.04705948230824
$ws.Cells.Item(6, 11).Value = 1.05237190659036
$ws.Cells.Item(6, 12).Value = 1.053313371240531
$ws.Cells.Item(6, 13).Value = 1.064580921135377
$ws.Cells.Item(6, 14).Value = 1.048546426872519
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.042366722221779
$ws.Cells.Item(7, 4).Value = 1.049747316515033
$ws.Cells.Item(7, 5).Value = 1.050661140631905
$ws.Cells.Item(7, 6).Value = 1.061944827980776
$ws.Cells.Item(7, 9).Value = 1.039035373821972
$ws.Cells.Item(7, 10).Value = 1.04686601492275
$ws.Cells.Item(7, 11).Value = 1.052196324577107
$ws.Cells.Item(7, 12).Value = 1.053107905158628
$ws.Cells.Item(7, 13).Value = 1.064364231773135
$ws.Cells.Item(7, 14).Value = 1.048352684741151
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.041092154745527
$ws.Cells.Item(8, 4).Value = 1.048765481381342
$ws.Cells.Item(8, 5).Value = 1.049555360175952
$ws.Cells.Item(8, 6).Value = 1.060794776778827
$ws.Cells.Item(8, 9).Value = 1.038817053059661
$ws.Cells.Item(8, 10).Value = 1.046057371534076
$ws.Cells.Item(8, 11).Value = 1.051461797362049
$ws.Cells.Item(8, 12).Value = 1.052249528474294
$ws.Cells.Item(8, 13).Value = 1.063458759795415
$ws.Cells.Item(8, 14).Value = 1.047542892986112
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.038851066359339
$ws.Cells.Item(9, 4).Value = 1.047038207271874
$ws.Cells.Item(9, 5).Value = 1.04761353404829
$ws.Cells.Item(9, 6).Value = 1.058774213910444
$ws.Cells.Item(9, 9).Value = 1.038423436869825
$ws.Cells.Item(9, 10).Value = 1.044632578062811
$ws.Cells.Item(9, 11).Value = 1.050165192507586
$ws.Cells.Item(9, 12).Value = 1.050738686655242
$ws.Cells.Item(9, 13).Value = 1.061864246221982
$ws.Cells.Item(9, 14).Value = 1.046116076144691
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.037360337171114
$ws.Cells.Item(10, 4).Value = 1.045888691656608
$ws.Cells.Item(10, 5).Value = 1.046323574199401
$ws.Cells.Item(10, 6).Value = 1.057431280721276
$ws.Cells.Item(10, 9).Value = 1.038155071259451
$ws.Cells.Item(10, 10).Value = 1.043682871272064
$ws.Cells.Item(10, 11).Value = 1.049299327206051
$ws.Cells.Item(10, 12).Value = 1.049732688866392
$ws.Cells.Item(10, 13).Value = 1.060802015682557
$ws.Cells.Item(10, 14).Value = 1.045165020661369
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.036715637132243
$ws.Cells.Item(11, 4).Value = 1.045391431541094
$ws.Cells.Item(11, 5).Value = 1.045766114212849
$ws.Cells.Item(11, 6).Value = 1.056850769577393
$ws.Cells.Item(11, 9).Value = 1.03803746062749
$ws.Cells.Item(11, 10).Value = 1.043271687459681
$ws.Cells.Item(11, 11).Value = 1.048924065129164
$ws.Cells.Item(11, 12).Value = 1.049297386768154
$ws.Cells.Item(11, 13).Value = 1.060342258801343
$ws.Cells.Item(11, 14).Value = 1.044753252920808
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.036476287303003
$ws.Cells.Item(12, 4).Value = 1.045206801731524
$ws.Cells.Item(12, 5).Value = 1.045559215619948
$ws.Cells.Item(12, 6).Value = 1.056635291827087
$ws.Cells.Item(12, 9).Value = 1.037993564025165
$ws.Cells.Item(12, 10).Value = 1.043118963324981
$ws.Cells.Item(12, 11).Value = 1.04878462659855
$ws.Cells.Item(12, 12).Value = 1.049135742625874
$ws.Cells.Item(12, 13).Value = 1.060171515332157
$ws.Cells.Item(12, 14).Value = 1.04460031190032
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.036527623187011
$ws.Cells.Item(13, 4).Value = 1.045246402046908
$ws.Cells.Item(13, 5).Value = 1.045603588496362
$ws.Cells.Item(13, 6).Value = 1.056681505726476
$ws.Cells.Item(13, 9).Value = 1.038002989524767
$ws.Cells.Item(13, 10).Value = 1.043151722799308
$ws.Cells.Item(13, 11).Value = 1.04881453885846
$ws.Cells.Item(13, 12).Value = 1.049170413723509
$ws.Cells.Item(13, 13).Value = 1.060208138975974
$ws.Cells.Item(13, 14).Value = 1.044633117896858
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.036695849937614
$ws.Cells.Item(14, 4).Value = 1.045376168450357
$ws.Cells.Item(14, 5).Value = 1.04574900849829
$ws.Cells.Item(14, 6).Value = 1.056832955051829
$ws.Cells.Item(14, 9).Value = 1.038033836417296
$ws.Cells.Item(14, 10).Value = 1.043259063069281
$ws.Cells.Item(14, 11).Value = 1.048912540099892
$ws.Cells.Item(14, 12).Value = 1.049284024259053
$ws.Cells.Item(14, 13).Value = 1.060328144459164
$ws.Cells.Item(14, 14).Value = 1.044740610602326
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.036799515982873
$ws.Cells.Item(15, 4).Value = 1.045456131761559
$ws.Cells.Item(15, 5).Value = 1.045838628715319
$ws.Cells.Item(15, 6).Value = 1.05692628789341
$ws.Cells.Item(15, 9).Value = 1.03805281429209
$ws.Cells.Item(15, 10).Value = 1.043325200021442
$ws.Cells.Item(15, 11).Value = 1.048972915374542
$ws.Cells.Item(15, 12).Value = 1.049354029647021
$ws.Cells.Item(15, 13).Value = 1.06040208790497
$ws.Cells.Item(15, 14).Value = 1.044806841476545
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.037403140631431
$ws.Cells.Item(16, 4).Value = 1.045921703573506
$ws.Cells.Item(16, 5).Value = 1.046360594304838
$ws.Cells.Item(16, 6).Value = 1.057469828276333
$ws.Cells.Item(16, 9).Value = 1.038162847101506
$ws.Cells.Item(16, 10).Value = 1.043710161231857
$ws.Cells.Item(16, 11).Value = 1.049324225114175
$ws.Cells.Item(16, 12).Value = 1.049761584883405
$ws.Cells.Item(16, 13).Value = 1.060832532469252
$ws.Cells.Item(16, 14).Value = 1.045192349376035
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.037781992082636
$ws.Cells.Item(17, 4).Value = 1.046213876082478
$ws.Cells.Item(17, 5).Value = 1.046688305347131
$ws.Cells.Item(17, 6).Value = 1.057811042338655
$ws.Cells.Item(17, 9).Value = 1.038231491583725
$ws.Cells.Item(17, 10).Value = 1.043951650118669
$ws.Cells.Item(17, 11).Value = 1.049544503065537
$ws.Cells.Item(17, 12).Value = 1.050017315055036
$ws.Cells.Item(17, 13).Value = 1.061102592328733
$ws.Cells.Item(17, 14).Value = 1.045434181204766
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038003046418217
$ws.Cells.Item(18, 4).Value = 1.046384342324806
$ws.Cells.Item(18, 5).Value = 1.046879559801622
$ws.Cells.Item(18, 6).Value = 1.058010161932888
$ws.Cells.Item(18, 9).Value = 1.038271395000576
$ws.Cells.Item(18, 10).Value = 1.044092510822737
$ws.Cells.Item(18, 11).Value = 1.049672954908685
$ws.Cells.Item(18, 12).Value = 1.050166507160114
$ws.Cells.Item(18, 13).Value = 1.06126013249718
$ws.Cells.Item(18, 14).Value = 1.04557524194719
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.038078433213867
$ws.Cells.Item(19, 4).Value = 1.046442474829128
$ws.Cells.Item(19, 5).Value = 1.046944790589088
$ws.Cells.Item(19, 6).Value = 1.058078072670568
$ws.Cells.Item(19, 9).Value = 1.038284977992818
$ws.Cells.Item(19, 10).Value = 1.044140541385171
$ws.Cells.Item(19, 11).Value = 1.049716748111772
$ws.Cells.Item(19, 12).Value = 1.050217382700864
$ws.Cells.Item(19, 13).Value = 1.061313852777315
$ws.Cells.Item(19, 14).Value = 1.045623340718532
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.037741336967647
$ws.Cells.Item(20, 4).Value = 1.04618252386939
$ws.Cells.Item(20, 5).Value = 1.046653134063975
$ws.Cells.Item(20, 6).Value = 1.057774423440455
$ws.Cells.Item(20, 9).Value = 1.038224140717592
$ws.Cells.Item(20, 10).Value = 1.043925740192209
$ws.Cells.Item(20, 11).Value = 1.049520872685848
$ws.Cells.Item(20, 12).Value = 1.049989874617713
$ws.Cells.Item(20, 13).Value = 1.06107361550567
$ws.Cells.Item(20, 14).Value = 1.045408234483238
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.036646308019268
$ws.Cells.Item(21, 4).Value = 1.045337953404016
$ws.Cells.Item(21, 5).Value = 1.045706181314362
$ws.Cells.Item(21, 6).Value = 1.056788352849021
$ws.Cells.Item(21, 9).Value = 1.038024758590405
$ws.Cells.Item(21, 10).Value = 1.043227453810616
$ws.Cells.Item(21, 11).Value = 1.048883682534169
$ws.Cells.Item(21, 12).Value = 1.049250567496988
$ws.Cells.Item(21, 13).Value = 1.060292804977569
$ws.Cells.Item(21, 14).Value = 1.044708956454889
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.035958517337491
$ws.Cells.Item(22, 4).Value = 1.044807371564969
$ws.Cells.Item(22, 5).Value = 1.045111760687091
$ws.Cells.Item(22, 6).Value = 1.05616923937211
$ws.Cells.Item(22, 9).Value = 1.037898179912173
$ws.Cells.Item(22, 10).Value = 1.042788459313526
$ws.Cells.Item(22, 11).Value = 1.048482770024634
$ws.Cells.Item(22, 12).Value = 1.048786005007813
$ws.Cells.Item(22, 13).Value = 1.05980205669872
$ws.Cells.Item(22, 14).Value = 1.044269338535259
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.036323061759536
$ws.Cells.Item(23, 4).Value = 1.045088601492755
$ws.Cells.Item(23, 5).Value = 1.045426782249821
$ws.Cells.Item(23, 6).Value = 1.056497360329415
$ws.Cells.Item(23, 9).Value = 1.037965397051292
$ws.Cells.Item(23, 10).Value = 1.043021173918551
$ws.Cells.Item(23, 11).Value = 1.048695328052768
$ws.Cells.Item(23, 12).Value = 1.04903225247484
$ws.Cells.Item(23, 13).Value = 1.060062194323211
$ws.Cells.Item(23, 14).Value = 1.044502383621716
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.037759707027217
$ws.Cells.Item(24, 4).Value = 1.046196690439999
$ws.Cells.Item(24, 5).Value = 1.046669026125075
$ws.Cells.Item(24, 6).Value = 1.057790969650305
$ws.Cells.Item(24, 9).Value = 1.038227462677074
$ws.Cells.Item(24, 10).Value = 1.043937447759434
$ws.Cells.Item(24, 11).Value = 1.04953155033798
$ws.Cells.Item(24, 12).Value = 1.050002273680616
$ws.Cells.Item(24, 13).Value = 1.061086708826692
$ws.Cells.Item(24, 14).Value = 1.045419958676551
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.039429858991272
$ws.Cells.Item(25, 4).Value = 1.047484402860384
$ws.Cells.Item(25, 5).Value = 1.048114739766811
$ws.Cells.Item(25, 6).Value = 1.059295860388297
$ws.Cells.Item(25, 9).Value = 1.038526248149303
$ws.Cells.Item(25, 10).Value = 1.045000898440454
$ws.Cells.Item(25, 11).Value = 1.050500658540748
$ws.Cells.Item(25, 12).Value = 1.051129063691954
$ws.Cells.Item(25, 13).Value = 1.062276334606261
$ws.Cells.Item(25, 14).Value = 1.046484919579518
